$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new value. NumberFormat is forced to
# text ("@") before assignment so numeric-looking strings (prices like
# "309.00" or "45.319.71") are preserved verbatim as text, matching the
# source data which stores these as inline strings, not numbers.
$updates = @(
    @{ Cell = "D2"; Value = "45.319.71" },
    @{ Cell = "E2"; Value = "  +6.03%  " },
    @{ Cell = "D3"; Value = "2.363.25" },
    @{ Cell = "E3"; Value = "  +2.19%  " },
    @{ Cell = "E4"; Value = "  +0.18%  " },
    @{ Cell = "D5"; Value = "109.59" },
    @{ Cell = "E5"; Value = "  +1.90%  " },
    @{ Cell = "D6"; Value = "309.00" },
    @{ Cell = "E6"; Value = "  -0.98%  " },
    @{ Cell = "E7"; Value = "  +0.24%  " },
    @{ Cell = "E8"; Value = "  -0.29%  " },
    @{ Cell = "E9"; Value = "  +1.46%  " },
    @{ Cell = "D10"; Value = "41.24" },
    @{ Cell = "E10"; Value = "  +2.35%  " },
    @{ Cell = "D11"; Value = "0.0919" },
    @{ Cell = "E11"; Value = "  +0.43%  " },
    @{ Cell = "E12"; Value = "  +1.00%  " },
    @{ Cell = "E13"; Value = "  +1.53%  " },
    @{ Cell = "E14"; Value = "  -1.22%  " },
    @{ Cell = "D15"; Value = "2.722.10" },
    @{ Cell = "E15"; Value = "  +2.15%  " },
    @{ Cell = "E16"; Value = "  -0.20%  " },
    @{ Cell = "D17"; Value = "2.360.31" },
    @{ Cell = "E17"; Value = "  +2.26%  " },
    @{ Cell = "D18"; Value = "45.244.41" },
    @{ Cell = "E18"; Value = "  +5.36%  " },
    @{ Cell = "D19"; Value = "7.31" },
    @{ Cell = "E19"; Value = "  -2.14%  " },
    @{ Cell = "E20"; Value = "  +0.84%  " },
    @{ Cell = "D21"; Value = "13.51" },
    @{ Cell = "E21"; Value = "  +2.59%  " },
    @{ Cell = "D22"; Value = "73.28" },
    @{ Cell = "E22"; Value = "  -0.43%  " },
    @{ Cell = "D23"; Value = "3.44" },
    @{ Cell = "E23"; Value = "  -1.35%  " },
    @{ Cell = "D24"; Value = "258.89" },
    @{ Cell = "E24"; Value = "  -2.66%  " },
    @{ Cell = "E25"; Value = "  +2.32%  " },
    @{ Cell = "E26"; Value = "  -0.41%  " },
    @{ Cell = "D27"; Value = "11.13" },
    @{ Cell = "E27"; Value = "  +1.14%  " },
    @{ Cell = "D28"; Value = "7.36" },
    @{ Cell = "E28"; Value = "  -3.92%  " },
    @{ Cell = "D29"; Value = "2.37" },
    @{ Cell = "E29"; Value = "  +2.88%  " },
    @{ Cell = "D30"; Value = "38.53" },
    @{ Cell = "E30"; Value = "  +0.04%  " },
    @{ Cell = "B31"; Value = "Hedera" },
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar" },
    @{ Cell = "D31"; Value = "0.0968" },
    @{ Cell = "E31"; Value = "  +11.15%  " },
    @{ Cell = "B32"; Value = "EthereumClassic" },
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" },
    @{ Cell = "D32"; Value = "22.43" },
    @{ Cell = "E32"; Value = "  +0.14%  " },
    @{ Cell = "D33"; Value = "170.58" },
    @{ Cell = "E33"; Value = "  +2.79%  " },
    @{ Cell = "E34"; Value = "  +6.43%  " },
    @{ Cell = "E35"; Value = "  +0.54%  " },
    @{ Cell = "D36"; Value = "4.86" },
    @{ Cell = "E36"; Value = "  +4.68%  " },
    @{ Cell = "E37"; Value = "  +1.64%  " },
    @{ Cell = "B38"; Value = "LidoDAOToken" },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo" },
    @{ Cell = "D38"; Value = "2.96" },
    @{ Cell = "E38"; Value = "  +4.49%  " },
    @{ Cell = "B39"; Value = "NEARProtocol" },
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" },
    @{ Cell = "D39"; Value = "3.93" },
    @{ Cell = "E39"; Value = "  +7.57%  " },
    @{ Cell = "B40"; Value = "VeChain" },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" },
    @{ Cell = "D40"; Value = "0.0358" },
    @{ Cell = "E40"; Value = "  -0.03%  " },
    @{ Cell = "D41"; Value = "1.73" },
    @{ Cell = "E41"; Value = "  +8.79%  " },
    @{ Cell = "D42"; Value = "99.28" },
    @{ Cell = "E42"; Value = "  -5.05%  " },
    @{ Cell = "D43"; Value = "0.233" },
    @{ Cell = "D44"; Value = "69.95" },
    @{ Cell = "E44"; Value = "  -1.67%  " },
    @{ Cell = "B45"; Value = "Celestia" },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia" },
    @{ Cell = "D45"; Value = "12.84" },
    @{ Cell = "E45"; Value = "  +3.77%  " },
    @{ Cell = "B46"; Value = "FirstDigitalUSD" },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd" },
    @{ Cell = "D46"; Value = "1.00" },
    @{ Cell = "E46"; Value = "  -0.43%  " },
    @{ Cell = "D47"; Value = "82.22" },
    @{ Cell = "E47"; Value = "  +7.38%  " },
    @{ Cell = "D48"; Value = "112.99" },
    @{ Cell = "E48"; Value = "  +0.30%  " },
    @{ Cell = "D49"; Value = "9.28" },
    @{ Cell = "E49"; Value = "  +4.78%  " },
    @{ Cell = "D50"; Value = "5.52" },
    @{ Cell = "E50"; Value = "  +5.03%  " },
    @{ Cell = "D51"; Value = "1.656.15" },
    @{ Cell = "E51"; Value = "  -0.67%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
